$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 1967.8572
$ws.Cells.Item(40, 9).Value = 1805
$ws.Cells.Item(40, 10).Value = 2375
$ws.Cells.Item(40, 11).Value = 1805
$ws.Cells.Item(40, 12).Value = 2375
$ws.Cells.Item(40, 13).Value = -1630
$ws.Cells.Item(40, 14).Value = -2725

# row 43 (Leve Item ID 5472)
$ws.Cells.Item(43, 8).Value = 6049.6
$ws.Cells.Item(43, 9).Value = 3000
$ws.Cells.Item(43, 10).Value = 8082.6665
$ws.Cells.Item(43, 11).Value = 3000
$ws.Cells.Item(43, 12).Value = 8082.6665
$ws.Cells.Item(43, 13).Value = -2931
$ws.Cells.Item(43, 14).Value = -8220.666499999999

# row 80 (Leve Item ID 12605)
$ws.Cells.Item(80, 8).Value = 307.60715
$ws.Cells.Item(80, 9).Value = 290.9375
$ws.Cells.Item(80, 10).Value = 329.83334
$ws.Cells.Item(80, 11).Value = 872.8125
$ws.Cells.Item(80, 12).Value = 989.5000200000001
$ws.Cells.Item(80, 13).Value = 125.1875
$ws.Cells.Item(80, 14).Value = -2985.50002

# row 83 (Leve Item ID 12605)
$ws.Cells.Item(83, 8).Value = 307.60715
$ws.Cells.Item(83, 9).Value = 290.9375
$ws.Cells.Item(83, 10).Value = 329.83334
$ws.Cells.Item(83, 11).Value = 2618.4375
$ws.Cells.Item(83, 12).Value = 2968.50006
$ws.Cells.Item(83, 13).Value = 2373.5625
$ws.Cells.Item(83, 14).Value = -12952.50006

# row 96 (Leve Item ID 19894)
$ws.Cells.Item(96, 8).Value = 2584.25
$ws.Cells.Item(96, 10).Value = 2999
$ws.Cells.Item(96, 12).Value = 8997
$ws.Cells.Item(96, 14).Value = -11743

# row 103 (Leve Item ID 19909)
$ws.Cells.Item(103, 8).Value = 382.5
$ws.Cells.Item(103, 9).Value = 290
$ws.Cells.Item(103, 11).Value = 870
$ws.Cells.Item(103, 13).Value = -284

# row 115 (Leve Item ID 27957)
$ws.Cells.Item(115, 8).Value = 2665.1667
$ws.Cells.Item(115, 9).Value = 2665.1667
$ws.Cells.Item(115, 11).Value = 7995.500100000001
$ws.Cells.Item(115, 13).Value = -6428.500100000001

# row 118 (Leve Item ID 27958)
$ws.Cells.Item(118, 8).Value = 199
$ws.Cells.Item(118, 9).Value = 199
$ws.Cells.Item(118, 11).Value = 597
$ws.Cells.Item(118, 13).Value = 1060

$ws = $wb.Worksheets.Item("ARM")
# row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 1408.1791
$ws.Cells.Item(74, 9).Value = 1094.258
$ws.Cells.Item(74, 11).Value = 1094.258
$ws.Cells.Item(74, 13).Value = -220.258

# row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 1408.1791
$ws.Cells.Item(77, 9).Value = 1094.258
$ws.Cells.Item(77, 11).Value = 5471.29
$ws.Cells.Item(77, 13).Value = -1103.29

# row 114 (Leve Item ID 25968)
$ws.Cells.Item(114, 8).Value = 100000
$ws.Cells.Item(114, 10).Value = 100000
$ws.Cells.Item(114, 12).Value = 100000
$ws.Cells.Item(114, 14).Value = -108678

# row 124 (Leve Item ID 34252)
$ws.Cells.Item(124, 8).Value = 51685.4
$ws.Cells.Item(124, 10).Value = 51685.4
$ws.Cells.Item(124, 12).Value = 51685.4
$ws.Cells.Item(124, 14).Value = -61505.4

# row 135 (Leve Item ID 42016)
$ws.Cells.Item(135, 8).Value = 44999.5
$ws.Cells.Item(135, 10).Value = 44999.5
$ws.Cells.Item(135, 12).Value = 44999.5
$ws.Cells.Item(135, 14).Value = -55139.5

$ws = $wb.Worksheets.Item("BSM")
# row 20 (Leve Item ID 14149)
$ws.Cells.Item(20, 8).Value = 8717
$ws.Cells.Item(20, 9).Value = 8988.799999999999
$ws.Cells.Item(20, 11).Value = 8988.799999999999
$ws.Cells.Item(20, 13).Value = -8741.799999999999

# row 135 (Leve Item ID 41992)
$ws.Cells.Item(135, 8).Value = 45000
$ws.Cells.Item(135, 10).Value = 45000
$ws.Cells.Item(135, 12).Value = 45000
$ws.Cells.Item(135, 14).Value = -55140

$ws = $wb.Worksheets.Item("CRP")
# row 62 (Leve Item ID 12580)
$ws.Cells.Item(62, 8).Value = 61995
$ws.Cells.Item(62, 10).Value = 137061.67
$ws.Cells.Item(62, 12).Value = 137061.67
$ws.Cells.Item(62, 14).Value = -138309.67

# row 65 (Leve Item ID 12580)
$ws.Cells.Item(65, 8).Value = 61995
$ws.Cells.Item(65, 10).Value = 137061.67
$ws.Cells.Item(65, 12).Value = 685308.3500000001
$ws.Cells.Item(65, 14).Value = -691548.3500000001

# row 99 (Leve Item ID 36198)
$ws.Cells.Item(99, 8).Value = 14482.36
$ws.Cells.Item(99, 9).Value = 11965.9
$ws.Cells.Item(99, 11).Value = 11965.9
$ws.Cells.Item(99, 13).Value = -10467.9

# row 126 (Leve Item ID 36198)
$ws.Cells.Item(126, 8).Value = 14482.36
$ws.Cells.Item(126, 9).Value = 11965.9
$ws.Cells.Item(126, 11).Value = 35897.7
$ws.Cells.Item(126, 13).Value = -33427.7

# row 130 (Leve Item ID 34689)
$ws.Cells.Item(130, 8).Value = 100000
$ws.Cells.Item(130, 10).Value = 100000
$ws.Cells.Item(130, 12).Value = 100000
$ws.Cells.Item(130, 14).Value = -110040

$ws = $wb.Worksheets.Item("CUL")
# row 2 (Leve Item ID 4847)
$ws.Cells.Item(2, 8).Value = 62553.312
$ws.Cells.Item(2, 10).Value = 120
$ws.Cells.Item(2, 12).Value = 720
$ws.Cells.Item(2, 14).Value = -946

# row 116 (Leve Item ID 27866)
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).ClearContents()  # M116 removed (was -1632.9998)
$ws.Cells.Item(116, 14).ClearContents()  # N116 removed (was -9581)

# row 137 (Leve Item ID 44088)
$ws.Cells.Item(137, 8).Value = 4952.143
$ws.Cells.Item(137, 9).Value = 3330.9
$ws.Cells.Item(137, 11).Value = 9992.700000000001
$ws.Cells.Item(137, 13).Value = -4892.700000000001

$ws = $wb.Worksheets.Item("GSM")
# row 80 (Leve Item ID 12521)
$ws.Cells.Item(80, 8).Value = 5187.5
$ws.Cells.Item(80, 9).Value = 5533.3335
$ws.Cells.Item(80, 10).Value = 4980
$ws.Cells.Item(80, 11).Value = 5533.3335
$ws.Cells.Item(80, 12).Value = 4980
$ws.Cells.Item(80, 13).Value = -4535.3335
$ws.Cells.Item(80, 14).Value = -6976

# row 83 (Leve Item ID 12521)
$ws.Cells.Item(83, 8).Value = 5187.5
$ws.Cells.Item(83, 9).Value = 5533.3335
$ws.Cells.Item(83, 10).Value = 4980
$ws.Cells.Item(83, 11).Value = 27666.6675
$ws.Cells.Item(83, 12).Value = 24900
$ws.Cells.Item(83, 13).Value = -22674.6675
$ws.Cells.Item(83, 14).Value = -34884

$ws = $wb.Worksheets.Item("LTW")
# row 7 (Leve Item ID 36249)
$ws.Cells.Item(7, 8).Value = 2672.1667
$ws.Cells.Item(7, 9).Value = 2706.6
$ws.Cells.Item(7, 11).Value = 2706.6
$ws.Cells.Item(7, 13).Value = -2594.6

# row 46 (Leve Item ID 5282)
$ws.Cells.Item(46, 8).Value = 3257.8823
$ws.Cells.Item(46, 9).Value = 1065.6666
$ws.Cells.Item(46, 10).Value = 3727.6428
$ws.Cells.Item(46, 11).Value = 1065.6666
$ws.Cells.Item(46, 12).Value = 3727.6428
$ws.Cells.Item(46, 13).Value = -877.6666
$ws.Cells.Item(46, 14).Value = -4103.6428

# row 68 (Leve Item ID 12563)
$ws.Cells.Item(68, 8).Value = 1998.5
$ws.Cells.Item(68, 9).Value = 1998.5
$ws.Cells.Item(68, 11).Value = 1998.5
$ws.Cells.Item(68, 13).Value = -1249.5

# row 71 (Leve Item ID 12563)
$ws.Cells.Item(71, 8).Value = 1998.5
$ws.Cells.Item(71, 9).Value = 1998.5
$ws.Cells.Item(71, 11).Value = 9992.5
$ws.Cells.Item(71, 13).Value = -6248.5

# row 82 (Leve Item ID 12565)
$ws.Cells.Item(82, 8).Value = 2793
$ws.Cells.Item(82, 10).Value = 1873.5555
$ws.Cells.Item(82, 12).Value = 1873.5555
$ws.Cells.Item(82, 14).Value = -2595.5555

# row 85 (Leve Item ID 12565)
$ws.Cells.Item(85, 8).Value = 2793
$ws.Cells.Item(85, 10).Value = 1873.5555
$ws.Cells.Item(85, 12).Value = 1873.5555
$ws.Cells.Item(85, 14).Value = -4369.5555

# row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 9723.75
$ws.Cells.Item(122, 9).Value = 9970
$ws.Cells.Item(122, 10).Value = 8000
$ws.Cells.Item(122, 11).Value = 29910
$ws.Cells.Item(122, 12).Value = 24000
$ws.Cells.Item(122, 13).Value = -27460
$ws.Cells.Item(122, 14).Value = -28900

# row 126 (Leve Item ID 36249)
$ws.Cells.Item(126, 8).Value = 2672.1667
$ws.Cells.Item(126, 9).Value = 2706.6
$ws.Cells.Item(126, 11).Value = 8119.799999999999
$ws.Cells.Item(126, 13).Value = -5649.799999999999

# row 141 (Leve Item ID 42487)
$ws.Cells.Item(141, 8).Value = 79571.336
$ws.Cells.Item(141, 10).Value = 84857
$ws.Cells.Item(141, 12).Value = 84857
$ws.Cells.Item(141, 14).Value = -95217

$ws = $wb.Worksheets.Item("WVR")
# row 81 (Leve Item ID 12596)
$ws.Cells.Item(81, 8).Value = 3071.8
$ws.Cells.Item(81, 9).Value = 3071.8
$ws.Cells.Item(81, 11).Value = 6143.6
$ws.Cells.Item(81, 13).Value = -5082.6

# row 84 (Leve Item ID 12596)
$ws.Cells.Item(84, 8).Value = 3071.8
$ws.Cells.Item(84, 9).Value = 3071.8
$ws.Cells.Item(84, 11).Value = 30718
$ws.Cells.Item(84, 13).Value = -25414

# row 126 (Leve Item ID 36210)
$ws.Cells.Item(126, 8).Value = 2161.5454
$ws.Cells.Item(126, 9).Value = 1623.8572
$ws.Cells.Item(126, 11).Value = 4871.571599999999
$ws.Cells.Item(126, 13).Value = -2401.571599999999

# row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 1564.683
$ws.Cells.Item(136, 9).Value = 1193.5135
$ws.Cells.Item(136, 11).Value = 3580.5405
$ws.Cells.Item(136, 13).Value = -1030.5405
